# Insert a new row at position 78 (pushes existing rows 78..121 down to 79..122,
# growing the used range from A1:R121 to A1:R122), then populate the new row
# with the latest weekly price-report entry for this market record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(78).Insert()

$ws.Cells.Item(78, 1).Value = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44981
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 100112030
$ws.Cells.Item(78, 7).Value = "Poroto granado"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 400
$ws.Cells.Item(78, 11).Value = 37000
$ws.Cells.Item(78, 12).Value = 38000
$ws.Cells.Item(78, 13).Value = 37500
$ws.Cells.Item(78, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(78, 16).Value = 1500
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
